# Generate Report for Handback
# Adds two new handed-back files to the report:
#   fd49e53a-52b8-4cbb-801e-53d992b87eec
#   17a28d78-1038-413d-a778-41bd32b30820
# as new rows 15/16 on the "Overview", "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

# Hyperlink-style font used throughout the workbook for "link" cells
# (underline + custom blue FF6495ED, matches the existing "HyperLink" cell style).
$hyperlinkColor = 15570276   # BGR packing of RGB(0x64,0x95,0xED)

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Data for the two new files being handed back.
# ---------------------------------------------------------------------------

$file1 = "fd49e53a-52b8-4cbb-801e-53d992b87eec"
$file1Hash = "3ec6b184af1efdb2ade0e50772c3fe0e903a2fbb"
$file1HandoffZh = "2016-03-02 15:33:21"
$file1HandbackZh = "2016-03-02 15:34:18"
$file1HandoffDe = "2016-03-02 15:33:34"
$file1HandbackDe = "2016-03-02 15:34:36"

$file2 = "17a28d78-1038-413d-a778-41bd32b30820"
$file2Hash = "2d3ee649a09fe4bdfd404add6bfc5ad9617f9bf3"

$statusInSync = "Handed back: in sync with en-US"
$reasonInclude = "Include"
$commit = "0000000000000000000000000000000000000000"

# ===========================================================================
# Sheet "Overview" (sheet 1) — columns: File Name | zh-cn | de-de
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 15 — fd49e53a...
$wsOverview.Hyperlinks.Add($wsOverview.Range("A15"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$file1.md", "", "", "$file1.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("A15")
$wsOverview.Range("B15").Value = $statusInSync
$wsOverview.Range("C15").Value = $statusInSync

# Row 16 — 17a28d78...
$wsOverview.Hyperlinks.Add($wsOverview.Range("A16"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$file2.md", "", "", "$file2.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("A16")
$wsOverview.Range("B16").Value = $statusInSync
$wsOverview.Range("C16").Value = $statusInSync

# ===========================================================================
# Sheet "zh-cn" (sheet 2) — columns:
# A Source File Name | B Status | C Correspond Handoff File |
# D Correspond Handoff Datetime | E Target File | F Correspond Handback File |
# G Correspond Handback DateTime | H Handoff Reason | I Dependency From
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 15 — fd49e53a...
$wsZh.Hyperlinks.Add($wsZh.Range("A15"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$file1.md", "", "", "$file1.md") | Out-Null
Style-AsHyperlink $wsZh.Range("A15")
$wsZh.Range("B15").Value = $statusInSync
$wsZh.Hyperlinks.Add($wsZh.Range("C15"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$file1.$file1Hash.zh-cn.xlf", "", "", "$file1.$file1Hash.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("C15")
$wsZh.Range("D15").Value = $file1HandoffZh
Style-AsDate $wsZh.Range("D15")
$wsZh.Hyperlinks.Add($wsZh.Range("E15"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$commit/e2e/$file1.md", "", "", "$file1.md") | Out-Null
Style-AsHyperlink $wsZh.Range("E15")
$wsZh.Hyperlinks.Add($wsZh.Range("F15"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$file1.$file1Hash.zh-cn.xlf", "", "", "$file1.$file1Hash.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("F15")
$wsZh.Range("G15").Value = $file1HandbackZh
$wsZh.Range("H15").Value = $reasonInclude

# Row 16 — 17a28d78...
$wsZh.Hyperlinks.Add($wsZh.Range("A16"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$file2.md", "", "", "$file2.md") | Out-Null
Style-AsHyperlink $wsZh.Range("A16")
$wsZh.Range("B16").Value = $statusInSync
$wsZh.Hyperlinks.Add($wsZh.Range("C16"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$file2.$file2Hash.zh-cn.xlf", "", "", "$file2.$file2Hash.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("C16")
$wsZh.Range("D16").Value = $file1HandoffZh
Style-AsDate $wsZh.Range("D16")
$wsZh.Hyperlinks.Add($wsZh.Range("E16"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$commit/e2e/$file2.md", "", "", "$file2.md") | Out-Null
Style-AsHyperlink $wsZh.Range("E16")
$wsZh.Hyperlinks.Add($wsZh.Range("F16"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$file2.$file2Hash.zh-cn.xlf", "", "", "$file2.$file2Hash.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("F16")
$wsZh.Range("G16").Value = $file1HandbackZh
$wsZh.Range("H16").Value = $reasonInclude

# ===========================================================================
# Sheet "de-de" (sheet 3) — same layout as "zh-cn"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

# Row 15 — fd49e53a...
$wsDe.Hyperlinks.Add($wsDe.Range("A15"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$file1.md", "", "", "$file1.md") | Out-Null
Style-AsHyperlink $wsDe.Range("A15")
$wsDe.Range("B15").Value = $statusInSync
$wsDe.Hyperlinks.Add($wsDe.Range("C15"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$file1.$file1Hash.de-de.xlf", "", "", "$file1.$file1Hash.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("C15")
$wsDe.Range("D15").Value = $file1HandoffDe
Style-AsDate $wsDe.Range("D15")
$wsDe.Hyperlinks.Add($wsDe.Range("E15"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$commit/e2e/$file1.md", "", "", "$file1.md") | Out-Null
Style-AsHyperlink $wsDe.Range("E15")
$wsDe.Hyperlinks.Add($wsDe.Range("F15"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$file1.$file1Hash.de-de.xlf", "", "", "$file1.$file1Hash.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("F15")
$wsDe.Range("G15").Value = $file1HandbackDe
$wsDe.Range("H15").Value = $reasonInclude

# Row 16 — 17a28d78...
$wsDe.Hyperlinks.Add($wsDe.Range("A16"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$file2.md", "", "", "$file2.md") | Out-Null
Style-AsHyperlink $wsDe.Range("A16")
$wsDe.Range("B16").Value = $statusInSync
$wsDe.Hyperlinks.Add($wsDe.Range("C16"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$file2.$file2Hash.de-de.xlf", "", "", "$file2.$file2Hash.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("C16")
$wsDe.Range("D16").Value = $file1HandoffDe
Style-AsDate $wsDe.Range("D16")
$wsDe.Hyperlinks.Add($wsDe.Range("E16"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$commit/e2e/$file2.md", "", "", "$file2.md") | Out-Null
Style-AsHyperlink $wsDe.Range("E16")
$wsDe.Hyperlinks.Add($wsDe.Range("F16"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$file2.$file2Hash.de-de.xlf", "", "", "$file2.$file2Hash.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("F16")
$wsDe.Range("G16").Value = $file1HandbackDe
$wsDe.Range("H16").Value = $reasonInclude

Write-Host "Report rows added for $file1 and $file2"
